$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")
$ws.Range("C2:C149").Value = 46075
